$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column D entirely (the "2024-01-09" column)
$ws.Range("D1:D19").EntireColumn.Delete()

# Update C column values for existing rows (2-19)
$ws.Range("C2").Value = 466633
$ws.Range("C3").Value = 312266
$ws.Range("C4").Value = 247702
$ws.Range("C5").Value = 57154
$ws.Range("C6").Value = 73873
$ws.Range("C7").Value = 389503
$ws.Range("C8").Value = 48937
$ws.Range("C9").Value = 132462
$ws.Range("C10").Value = 95898
$ws.Range("C11").Value = 146109
$ws.Range("C12").Value = 52695
$ws.Range("C13").Value = 3980296
$ws.Range("C14").Value = 79928
$ws.Range("C15").Value = 2231638
$ws.Range("C16").Value = 14154533
$ws.Range("C17").Value = 36
$ws.Range("C18").Value = 4323681
$ws.Range("C19").Value = 5270698

# Update B6 name
$ws.Range("B6").Value = "Giovanna Pitel"

# Add new rows 20-27
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "Davi"
$ws.Range("C20").Value = 292415

$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "Giovanna"
$ws.Range("C21").Value = 201669

$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "Isabelle"
$ws.Range("C22").Value = 614226

$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "Juninho"
$ws.Range("C23").Value = 21799

$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "Lucas Henrique"
$ws.Range("C24").Value = 38615

$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "Michel"
$ws.Range("C25").Value = 25078

$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "Raqueli"
$ws.Range("C26").Value = 90877

$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "Thalyta"
$ws.Range("C27").Value = 29370
